# Auto-generated edit script: rolls quarterly data forward by one quarter
# (drops oldest quarter column, shifts remaining left, appends new quarter data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update quarter header labels (row 8, 30, 52, 68, 90), columns E..N ---
$quarterLabels = @(
    "فصل چهارم منتهی به 1399/09",
    "فصل اول منتهی به 1399/12",
    "فصل دوم منتهی به 1400/03",
    "فصل سوم منتهی به 1400/06",
    "فصل چهارم منتهی به 1400/09",
    "فصل اول منتهی به 1400/12",
    "فصل دوم منتهی به 1401/03",
    "فصل سوم منتهی به 1401/06",
    "فصل چهارم منتهی به 1401/09",
    "فصل اول منتهی به 1401/12"
)

$headerRows = @(8, 30, 52, 68, 90)
foreach ($hr in $headerRows) {
    $col = 5
    foreach ($label in $quarterLabels) {
        $ws.Cells.Item($hr, $col).Value = $label
        $col = $col + 1
    }
}

# --- Update data rows, columns E..N ---
$rowData = @{
    11 = @(9159, 6450, 5349, 6179, 9282, 3507, 2809, 3818, 3503, 2427)
    12 = @(14469, 9621, 5559, 12607, 15245, 5760, 7215, 6154, 8125, 7192)
    13 = @(6950, 6336, 4224, 6501, 7323, 3461, 3741, 4341, 4503, 4598)
    14 = @(30578, 22407, 15132, 25287, 31850, 12728, 13765, 14313, 16131, 14217)
    16 = @(2865, 1318, 2300, 1922, 2087, 2621, 1134, 1622, 1380, 1289)
    17 = @(5445, 5382, 6372, 5766, 6628, 3636, 7443, 14119, 8995, 6706)
    18 = @(144, 23, 18, 103, 349, 232, 352, 643, 61, 63)
    19 = @(8454, 6723, 8690, 7791, 9064, 6489, 8929, 16384, 10436, 8058)
    21 = @(-5470, 0, -2206, 1478, 1359, 845, 519, 608, "-", "-")
    22 = @("-", "-", "-", "-", "-", "-", "-", "-", 2759, 2301)
    23 = @(-5470, 0, -2206, 1478, 1359, 845, 519, 608, 2759, 2301)
    24 = @(0, 0, 0, "-", "-", "-", "-", "-", "-", "-")
    25 = @(0, 0, 0, 0, "-", "-", 0, 0, 0, 0)
    26 = @(33562, 29130, 21616, 34556, 42273, 20062, 23213, 31305, 29326, 24576)
    33 = @(1374776, 1488707, 1372659, 1743708, 2799061, 1108636, 1077066, 1739842, 1269607, 879191)
    34 = @(1743038, 1318369, 744059, 2055009, 2365801, 956709, 1236229, 1389603, 1778117, 1591750)
    35 = @(1788522, 1818110, 1221209, 2205991, 2546224, 1257623, 1436935, 1916784, 2182135, 2434111)
    36 = @(4906336, 4625186, 3337927, 6004708, 7711086, 3322968, 3750230, 5046229, 5229859, 4905052)
    38 = @(459674, 290586, 434829, 500686, 338102, 688098, 299016, 409992, 189515, 397192)
    39 = @(353038, 284993, 379963, 380596, 467146, 242300, 417570, 1177608, 878380, 894044)
    40 = @(26141, 4036, 4280, 26291, 100092, 78311, 131473, 284529, -10319, 32017)
    41 = @(838853, 579615, 819072, 907573, 905340, 1008709, 848059, 1872129, 1057576, 1323253)
    43 = @(44784, 50314, 39687, 39300, 36241, 25449, 17895, 21059, "-", "-")
    44 = @("-", "-", "-", "-", "-", "-", "-", "-", 130621, 129789)
    45 = @(44784, 50314, 39687, 39300, 36241, 25449, 17895, 21059, 130621, 129789)
    46 = @(0, 0, 0, "-", "-", "-", "-", "-", "-", "-")
    47 = @(0, 0, 0, 0, "-", "-", 0, 0, 0, 0)
    48 = @(5789973, 5255115, 4196686, 6951581, 8652667, 4357126, 4616184, 6939417, 6418056, 6358094)
    55 = @(150101103, 230807287, 256619742, 282199061, 301557962, 316120901, 383433962, 455694605, 362434199, 362254223)
    56 = @(120467068, 137030350, 133847634, 163005394, 155185372, 166095313, 171341511, 225804842, 218845169, 221322303)
    57 = @(257341295, 286949179, 289111979, 339331026, 347702308, 363369835, 384104518, 441553559, 484595825, 529384732)
    59 = @(160444677, 220474962, 189056087, 260502601, 162003833, 262532621, 263682540, 252769420, 137329710, 308139643)
    60 = @(64837098, 52952991, 59630100, 66006937, 70480688, 66639164, 56102378, 83405907, 97652029, 133320012)
    61 = @(181534722, 175478261, 237777778, 255252427, 286796562, 337547414, 373502841, 442502333, 169163934, 508206349)
    63 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    64 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", 56405476)
    71 = @(-1003436, -1253516, -1134794, -1455743, -2267023, -901403, -841235, -1324007, -1054995, -787475)
    72 = @(-1394462, -1152270, -732130, -1820297, -2369290, -1000628, -1334710, -1468696, -1780783, -1617819)
    73 = @(-1437497, -1669883, -1132778, -1770048, -2301732, -1244183, -1526195, -1921940, -1979529, -2414947)
    74 = @(-3835395, -4075669, -2999702, -5046088, -6938045, -3146214, -3702140, -4714643, -4815307, -4820241)
    76 = @(-339861, -260701, -433799, -426259, -294036, -697035, -287497, -353771, -243212, -306304)
    77 = @(-177337, -226092, -330624, -313156, -399109, -232535, -576070, -1146257, -757789, -635190)
    78 = @(-15155, -4844, -4311, -23633, -82127, -85463, -110164, -291295, -14578, -19538)
    79 = @(-532353, -491637, -768734, -763048, -775272, -1015033, -973731, -1791323, -1015579, -961032)
    81 = @(-11383, -13985, -25269, 0, 0, -14020, -9735, -25722, "-", "-")
    82 = @("-", "-", "-", "-", "-", "-", "-", "-", -67424, -84044)
    83 = @(-11383, -13985, -25269, 0, 0, -14020, -9735, -25722, -67424, -84044)
    84 = @(0, 0, 0, "-", "-", "-", "-", "-", "-", "-")
    85 = @(0, 0, 0, 0, "-", "-", 0, 0, 0, 0)
    86 = @(-4379131, -4581291, -3793705, -5809136, -7713317, -4175267, -4685606, -6531688, -5898310, -5865317)
    93 = @(371340, 235191, 237865, 287965, 532038, 207233, 235831, 415835, 214612, 91716)
    94 = @(348576, 166099, 11929, 234712, -3489, -43919, -98481, -79093, -2666, -26069)
    95 = @(351025, 148227, 88431, 435943, 244492, 13440, -89260, -5156, 202606, 19164)
    96 = @(1070941, 549517, 338225, 958620, 773041, 176754, 48090, 331586, 414552, 84811)
    98 = @(119813, 29885, 1030, 74427, 44066, -8937, 11519, 56221, -53697, 90888)
    99 = @(175701, 58901, 49339, 67440, 68037, 9765, -158500, 31351, 120591, 258854)
    100 = @(10986, -808, -31, 2658, 17965, -7152, 21309, -6766, -24897, 12479)
    101 = @(306500, 87978, 50338, 144525, 130068, -6324, -125672, 80806, 41997, 362221)
    103 = @(33401, 36329, 14418, 39300, 36241, 11429, 8160, -4663, "-", "-")
    104 = @("-", "-", "-", "-", "-", "-", "-", "-", 63197, 45745)
    105 = @(33401, 36329, 14418, 39300, 36241, 11429, 8160, -4663, 63197, 45745)
    106 = @(1410842, 673824, 402981, 1142445, 939350, 181859, -69422, 407729, 519746, 492777)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $col = 5
    foreach ($v in $vals) {
        $ws.Cells.Item([int]$r, $col).Value = $v
        $col = $col + 1
    }
}
